$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (price) to be treated as text so values like "0.0000280"
# or "3.10" keep their exact formatting instead of being auto-converted to numbers.
$ws.Range("D2:D50").NumberFormat = "@"

$ws.Range("D2").Value = '70.808.96'
$ws.Range("E2").Value = '  +2.61%  '

$ws.Range("D3").Value = '3.595.59'
$ws.Range("E3").Value = '  +2.34%  '

$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").Value = '604.97'
$ws.Range("E5").Value = '  +2.76%  '

$ws.Range("D6").Value = '174.39'
$ws.Range("E6").Value = '  +1.47%  '

$ws.Range("D7").Value = '3.592.21'
$ws.Range("E7").Value = '  +2.54%  '

$ws.Range("E8").Value = '  +1.10%  '

$ws.Range("E9").Value = '  -0.05%  '

$ws.Range("D10").Value = '0.202'
$ws.Range("E10").Value = '  +6.94%  '

$ws.Range("D11").Value = '7.47'
$ws.Range("E11").Value = '  +9.57%  '

$ws.Range("D12").Value = '0.592'
$ws.Range("E12").Value = '  +1.79%  '

$ws.Range("D13").Value = '47.24'
$ws.Range("E13").Value = '  -0.16%  '

$ws.Range("D14").Value = '0.0000280'
$ws.Range("E14").Value = '  +1.89%  '

$ws.Range("D15").Value = '4.172.53'
$ws.Range("E15").Value = '  +2.28%  '

$ws.Range("D16").Value = '8.47'
$ws.Range("E16").Value = '  -0.23%  '

$ws.Range("D17").Value = '618.66'
$ws.Range("E17").Value = '  -1.37%  '

$ws.Range("D18").Value = '3.594.69'
$ws.Range("E18").Value = '  +2.69%  '

$ws.Range("D19").Value = '70.938.42'
$ws.Range("E19").Value = '  +2.63%  '

$ws.Range("E20").Value = '  -2.13%  '

$ws.Range("D21").Value = '17.55'
$ws.Range("E21").Value = '  +0.85%  '

$ws.Range("D22").Value = '0.893'
$ws.Range("E22").Value = '  +0.66%  '

$ws.Range("D23").Value = '9.27'
$ws.Range("E23").Value = '  -16.81%  '

$ws.Range("D24").Value = '16.16'
$ws.Range("E24").Value = '  +1.52%  '

$ws.Range("D25").Value = '98.02'
$ws.Range("E25").Value = '  +1.00%  '

$ws.Range("D26").Value = '3.82'
$ws.Range("E26").Value = '  +0.15%  '

$ws.Range("B27").Value = 'ImmutableX'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D27").Value = '2.67'
$ws.Range("E27").Value = '  +1.53%  '

$ws.Range("B28").Value = 'Dai'
$ws.Range("C28").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D28").Value = '0.999'
$ws.Range("E28").Value = '  -0.05%  '

$ws.Range("B29").Value = 'EthereumClassic'
$ws.Range("C29").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D29").Value = '34.21'
$ws.Range("E29").Value = '  +4.62%  '

$ws.Range("B30").Value = 'RenderToken'
$ws.Range("C30").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D30").Value = '9.35'
$ws.Range("E30").Value = '  +1.01%  '

$ws.Range("D31").Value = '8.56'
$ws.Range("E31").Value = '  +0.07%  '

$ws.Range("D32").Value = '3.10'
$ws.Range("E32").Value = '  -1.44%  '

$ws.Range("D33").Value = '7.26'
$ws.Range("E33").Value = '  +4.71%  '

$ws.Range("E34").Value = '  -1.48%  '

$ws.Range("D35").Value = '628.25'
$ws.Range("E35").Value = '  -1.51%  '

$ws.Range("D36").Value = '3.77'
$ws.Range("E36").Value = '  +8.04%  '

$ws.Range("D37").Value = '0.103'
$ws.Range("E37").Value = '  -0.15%  '

$ws.Range("D38").Value = '10.91'
$ws.Range("E38").Value = '  +1.41%  '

$ws.Range("E39").Value = '  +7.50%  '

$ws.Range("D40").Value = '57.66'
$ws.Range("E40").Value = '  +0.81%  '

$ws.Range("E41").Value = '  +0.05%  '

$ws.Range("D42").Value = '0.145'
$ws.Range("E42").Value = '  +7.15%  '

$ws.Range("D43").Value = '3.406.47'
$ws.Range("E43").Value = '  +0.51%  '

$ws.Range("D44").Value = '0.326'
$ws.Range("E44").Value = '  -0.47%  '

$ws.Range("B45").Value = 'PEPE'
$ws.Range("C45").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D45").Value = '0.0₃0721'
$ws.Range("E45").Value = '  +3.19%  '

$ws.Range("B46").Value = 'ThetaToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D46").Value = '3.02'
$ws.Range("E46").Value = '  +10.42%  '

$ws.Range("E47").Value = '  +6.77%  '

$ws.Range("D48").Value = '33.11'
$ws.Range("E48").Value = '  +0.97%  '

$ws.Range("E49").Value = '  +0.99%  '

$ws.Range("D50").Value = '132.83'
$ws.Range("E50").Value = '  +0.22%  '

# Reset the number format back to General/Normal style so the saved file
# does not carry a spurious text-format style on column D.
$ws.Range("D2:D50").Style = "Normal"
